$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename genus "Erythroparvovirus" -> "Erythyroparvovirus" for existing rows (2-8)
$ws.Range("F2:F8").Value = "Erythyroparvovirus"

# Add new row 9 data (Ungulate erythroparvovirus 2 / GiEPV / ErthryoPV-Giraffe-1)
$ws.Range("A9").Value = "ErthryoPV-Giraffe-1"
$ws.Range("B9").Value = "GiEPV"
$ws.Range("C9").Value = "Ungulate erythroparvovirus 2"
$ws.Range("D9").Value = "NK"
$ws.Range("E9").Value = "NK"
$ws.Range("F9").Value = "Erythyroparvovirus"
$ws.Range("G9").Value = "NK"
$ws.Range("H9").Value = "NK"
$ws.Range("I9").Value = "NK"
$ws.Range("J9").Value = "NK"
$ws.Range("K9").Value = "NK"
$ws.Range("L9").Value = $false

$ws.Range("D21").Select() | Out-Null
